$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "last updated" timestamp in the title cell (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 27 de Marzo de 2020 a las 22:28"

# --- Swap Camerun / Estado de Palestina rows (106 <-> 107) and refresh their stats ---
$ws.Range("A106").Value = "Camerun"
$ws.Range("C106").Value = 16
$ws.Range("D106").Value = 2
$ws.Range("E106").Value = 87
$ws.Range("G106").Value = 1
$ws.Range("H106").Value = 2

$ws.Range("A107").Value = "Estado de Palestina"
$ws.Range("B107").Value = 91
$ws.Range("C107").Value = 5
$ws.Range("D107").Value = 17
$ws.Range("E107").Value = 73
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 1

# --- Update country case numbers (row 4: Estados Unidos) ---
$ws.Range("B4").Value = 100514
$ws.Range("C4").Value = 15079
$ws.Range("E4").Value = 96503
$ws.Range("G4").Value = 251
$ws.Range("H4").Value = 1546

# --- Row 7: Espana ---
$ws.Range("B7").Value = 64285
$ws.Range("C7").Value = 6499
$ws.Range("E7").Value = 49988
$ws.Range("G7").Value = 575
$ws.Range("H7").Value = 4940

# --- Row 23: Suecia ---
$ws.Range("B23").Value = 3069
$ws.Range("C23").Value = 229
$ws.Range("E23").Value = 2961

# --- Row 68: Marruecos ---
$ws.Range("B68").Value = 345
$ws.Range("C68").Value = 70
$ws.Range("E68").Value = 311
$ws.Range("G68").Value = 12
$ws.Range("H68").Value = 23

# --- Row 90: Republica de Chipre ---
$ws.Range("D90").Value = 15
$ws.Range("E90").Value = 142
